$d = $word.ActiveDocument

# --- Step 1: Remove the "Meta description" paragraph that currently follows the
#     title (Heading1) paragraph at the top of the document. ---
$metaText = "Meta description: Read our review of Da Vinci Diamonds and play it for free. Discover its tumbling reels, bonus features, and high-paying symbols with a bet up to €100."

$metaParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext.TrimEnd("`r") -eq $metaText) {
        $metaParaIndex = $i
        break
    }
}
if ($metaParaIndex -eq -1) {
    throw "Could not locate the 'Meta description' paragraph"
}

$metaPara = $d.Paragraphs.Item($metaParaIndex)
$metaRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End)
$metaRange.Delete()

# --- Step 2: Replace the "Create a feature image..." image-prompt paragraph
#     text (now the last paragraph) with the meta-description copy, keeping the
#     existing italic run formatting intact. ---
$oldText = "Create a feature image for Da Vinci Diamonds that incorporates a happy Maya warrior with glasses in a cartoon style. The Maya warrior should be smiling and holding up a diamond-shaped icon similar to the symbol in the game. The background should feature a museum-like setting with displays of Da Vinci's masterpieces, precious stones, and other slot machine symbols. The overall design should be colorful and eye-catching, with a mix of ancient and modern elements that represent the game's theme."
$newText = "Read our review of Da Vinci Diamonds and play it for free. Discover its tumbling reels, bonus features, and high-paying symbols with a bet up to €100."

$replaced = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $replaced) {
    throw "Could not find the image-prompt paragraph text to replace"
}

# --- Step 3: Insert a new paragraph with the bold title text right before that
#     final (now re-worded, italic) paragraph. We inject raw WordOpenXML so the
#     run layout matches exactly (leading empty run + single bold run), then
#     split it away from the following paragraph with InsertParagraphBefore. ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertStart = $lastPara.Range.Start
$insertRange = $d.Range($insertStart, $insertStart)

$titleText = "Play Da Vinci Diamonds Free: Check the Review | Max €100 per Line"

$xmlPayload = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>REPLACE_TITLE</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$xmlPayload = $xmlPayload.Replace("REPLACE_TITLE", $titleText)

$insertRange.InsertXML($xmlPayload)

# Split the merged paragraph right after the title text so the bold title and
# the (now re-worded) italic paragraph end up as two separate paragraphs.
$splitPos = $insertStart + $titleText.Length
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphBefore()

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
